# Auto-generated edit script: updates currentAveragePrice / profit columns (H-N)
# across all 8 Leve-crafting job sheets per scheduled price-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 785.2
$ws.Range("I2").Value = 1315.5
$ws.Range("J2").Value = 431.66666
$ws.Range("K2").Value = 1315.5
$ws.Range("L2").Value = 431.66666
$ws.Range("M2").Value = -1202.5
$ws.Range("N2").Value = -657.66666

# Row 4
$ws.Range("H4").Value = 398.7143
$ws.Range("I4").Value = 91.75
$ws.Range("J4").Value = 808
$ws.Range("K4").Value = 91.75
$ws.Range("L4").Value = 808
$ws.Range("M4").Value = 22.25
$ws.Range("N4").Value = -1036

# Row 9
$ws.Range("H9").Value = 463.89743
$ws.Range("I9").Value = 512.9091
$ws.Range("K9").Value = 512.9091
$ws.Range("M9").Value = -343.9091

# Row 12
$ws.Range("H12").Value = 1228.2307
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = 1322.25
$ws.Range("K12").Value = 100
$ws.Range("L12").Value = 1322.25
$ws.Range("M12").Value = 70
$ws.Range("N12").Value = -1662.25

# Row 64
$ws.Range("H64").Value = 3174.487
$ws.Range("I64").Value = 2874.5
$ws.Range("J64").Value = 3251.9033
$ws.Range("K64").Value = 2874.5
$ws.Range("L64").Value = 3251.9033
$ws.Range("M64").Value = -2626.5
$ws.Range("N64").Value = -3747.9033

# Row 67
$ws.Range("H67").Value = 3174.487
$ws.Range("I67").Value = 2874.5
$ws.Range("J67").Value = 3251.9033
$ws.Range("K67").Value = 2874.5
$ws.Range("L67").Value = 3251.9033
$ws.Range("M67").Value = -2016.5
$ws.Range("N67").Value = -4967.9033

# Row 74
$ws.Range("H74").Value = 4189.579
$ws.Range("I74").Value = 4085.4666
$ws.Range("J74").Value = 4257.478
$ws.Range("K74").Value = 4085.4666
$ws.Range("L74").Value = 4257.478
$ws.Range("M74").Value = -3149.4666
$ws.Range("N74").Value = -6129.478

# Row 76
$ws.Range("H76").Value = 188785.67
$ws.Range("I76").Value = 280253.5
$ws.Range("K76").Value = 280253.5
$ws.Range("M76").Value = -279938.5

# Row 77
$ws.Range("H77").Value = 4189.579
$ws.Range("I77").Value = 4085.4666
$ws.Range("J77").Value = 4257.478
$ws.Range("K77").Value = 20427.333
$ws.Range("L77").Value = 21287.39
$ws.Range("M77").Value = -15747.333
$ws.Range("N77").Value = -30647.39

# Row 79
$ws.Range("H79").Value = 188785.67
$ws.Range("I79").Value = 280253.5
$ws.Range("K79").Value = 280253.5
$ws.Range("M79").Value = -279161.5

# Row 129
$ws.Range("H129").Value = 1072.875
$ws.Range("J129").Value = 1335.9445
$ws.Range("L129").Value = 4007.8335
$ws.Range("N129").Value = -14007.8335

# Row 137
$ws.Range("H137").Value = 1552.5518
$ws.Range("I137").Value = 1264.8948
$ws.Range("J137").Value = 2099.1
$ws.Range("K137").Value = 3794.6844
$ws.Range("L137").Value = 6297.299999999999
$ws.Range("M137").Value = -1244.6844
$ws.Range("N137").Value = -11397.3

# Row 138
$ws.Range("H138").Value = 1573.3387
$ws.Range("I138").Value = 644.4
$ws.Range("K138").Value = 1933.2
$ws.Range("M138").Value = 3206.8

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12248.297
$ws.Range("I32").Value = 2741.6511
$ws.Range("K32").Value = 2741.6511
$ws.Range("M32").Value = -2454.6511

# Row 45
$ws.Range("H45").Value = 1698.7368
$ws.Range("I45").Value = 1832.4166
$ws.Range("J45").Value = 1469.5714
$ws.Range("K45").Value = 1832.4166
$ws.Range("L45").Value = 1469.5714
$ws.Range("M45").Value = -1455.4166
$ws.Range("N45").Value = -2223.5714

# Row 63
$ws.Range("H63").Value = 2708.3333
$ws.Range("I63").Value = 1890
$ws.Range("J63").Value = 6800
$ws.Range("K63").Value = 1890
$ws.Range("L63").Value = 6800
$ws.Range("M63").Value = -1204
$ws.Range("N63").Value = -8172

# Row 66
$ws.Range("H66").Value = 2708.3333
$ws.Range("I66").Value = 1890
$ws.Range("J66").Value = 6800
$ws.Range("K66").Value = 9450
$ws.Range("L66").Value = 34000
$ws.Range("M66").Value = -6018
$ws.Range("N66").Value = -40864

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2662.2222
$ws.Range("I105").Value = 2441.9048
$ws.Range("J105").Value = 3433.3333
$ws.Range("K105").Value = 2441.9048
$ws.Range("L105").Value = 3433.3333
$ws.Range("M105").Value = -694.9047999999998
$ws.Range("N105").Value = -6927.3333

# Row 108
$ws.Range("H108").Value = 40000
$ws.Range("J108").Value = 40000
$ws.Range("L108").Value = 40000
$ws.Range("N108").Value = -47680

$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 20000
$ws.Range("I23").Value = 5000
$ws.Range("J23").Value = 27500
$ws.Range("K23").Value = 5000
$ws.Range("L23").Value = 27500
$ws.Range("M23").Value = -4760
$ws.Range("N23").Value = -27980

# Row 27
$ws.Range("H27").Value = 20000
$ws.Range("I27").Value = 5000
$ws.Range("J27").Value = 27500
$ws.Range("K27").Value = 5000
$ws.Range("L27").Value = 27500
$ws.Range("M27").Value = -4808
$ws.Range("N27").Value = -27884

# Row 31
$ws.Range("H31").Value = 8912082
$ws.Range("I31").Value = 4786636
$ws.Range("J31").Value = 66668332
$ws.Range("K31").Value = 4786636
$ws.Range("L31").Value = 66668332
$ws.Range("M31").Value = -4786341
$ws.Range("N31").Value = -66668922

# Row 34
$ws.Range("H34").Value = 8912082
$ws.Range("I34").Value = 4786636
$ws.Range("J34").Value = 66668332
$ws.Range("K34").Value = 4786636
$ws.Range("L34").Value = 66668332
$ws.Range("M34").Value = -4786434
$ws.Range("N34").Value = -66668736

# Row 62
$ws.Range("H62").Value = 142859870
$ws.Range("I62").Value = 3147.5
$ws.Range("J62").Value = 333335500
$ws.Range("K62").Value = 3147.5
$ws.Range("L62").Value = 333335500
$ws.Range("M62").Value = -2523.5
$ws.Range("N62").Value = -333336748

# Row 65
$ws.Range("H65").Value = 142859870
$ws.Range("I65").Value = 3147.5
$ws.Range("J65").Value = 333335500
$ws.Range("K65").Value = 15737.5
$ws.Range("L65").Value = 1666677500
$ws.Range("M65").Value = -12617.5
$ws.Range("N65").Value = -1666683740

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 874.64
$ws.Range("J131").Value = 950.8372000000001
$ws.Range("L131").Value = 2852.5116
$ws.Range("N131").Value = -12932.5116

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5147.2573
$ws.Range("I70").Value = 5197.609
$ws.Range("J70").Value = 5050.75
$ws.Range("K70").Value = 5197.609
$ws.Range("L70").Value = 5050.75
$ws.Range("M70").Value = -4927.609
$ws.Range("N70").Value = -5590.75

# Row 73
$ws.Range("H73").Value = 5147.2573
$ws.Range("I73").Value = 5197.609
$ws.Range("J73").Value = 5050.75
$ws.Range("K73").Value = 5197.609
$ws.Range("L73").Value = 5050.75
$ws.Range("M73").Value = -4261.609
$ws.Range("N73").Value = -6922.75

# Row 80
$ws.Range("H80").Value = 3698.2273
$ws.Range("I80").Value = 4732.778
$ws.Range("J80").Value = 2982
$ws.Range("K80").Value = 4732.778
$ws.Range("L80").Value = 2982
$ws.Range("M80").Value = -3734.778
$ws.Range("N80").Value = -4978

# Row 83
$ws.Range("H83").Value = 3698.2273
$ws.Range("I83").Value = 4732.778
$ws.Range("J83").Value = 2982
$ws.Range("K83").Value = 23663.89
$ws.Range("L83").Value = 14910
$ws.Range("M83").Value = -18671.89
$ws.Range("N83").Value = -24894

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1201.591
$ws.Range("I16").Value = 751
$ws.Range("J16").Value = 2733.6
$ws.Range("K16").Value = 751
$ws.Range("L16").Value = 2733.6
$ws.Range("M16").Value = -581
$ws.Range("N16").Value = -3073.6

# Row 22
$ws.Range("H22").Value = 448.75
$ws.Range("I22").Value = 297.5
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 297.5
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = -2.5
$ws.Range("N22").Value = -1190

# Row 27
$ws.Range("H27").Value = 448.75
$ws.Range("I27").Value = 297.5
$ws.Range("J27").Value = 600
$ws.Range("K27").Value = 297.5
$ws.Range("L27").Value = 600
$ws.Range("M27").Value = -190.5
$ws.Range("N27").Value = -814

# Row 55
$ws.Range("H55").Value = 500
$ws.Range("I55").Value = 200
$ws.Range("K55").Value = 200
$ws.Range("M55").Value = -27

# Row 122
$ws.Range("H122").Value = 4978.706
$ws.Range("I122").Value = 3694.3635
$ws.Range("J122").Value = 7333.3335
$ws.Range("K122").Value = 11083.0905
$ws.Range("L122").Value = 22000.0005
$ws.Range("M122").Value = -8633.0905
$ws.Range("N122").Value = -26900.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 20077
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

# Row 81
$ws.Range("H81").Value = 58825516
$ws.Range("I81").Value = 66668452
$ws.Range("J81").Value = 3499.5
$ws.Range("K81").Value = 133336904
$ws.Range("L81").Value = 6999
$ws.Range("M81").Value = -133335843
$ws.Range("N81").Value = -9121

# Row 84
$ws.Range("H84").Value = 58825516
$ws.Range("I84").Value = 66668452
$ws.Range("J84").Value = 3499.5
$ws.Range("K84").Value = 666684520
$ws.Range("L84").Value = 34995
$ws.Range("M84").Value = -666679216
$ws.Range("N84").Value = -45603
